# Mifos automation workbook correction
# - Transactions sheet: refresh the transaction listing with the latest
#   accrual/repayment/disbursement rows (9 rows instead of 6), drop the
#   now-unused K/L helper columns, and move the selection.
# - DeleteAccountClosures sheet: remove the stray "Navigate"/"LoanAccount"
#   helper row that was left over, and make this the active sheet.
# - Repay1 sheet: selection cursor nudged up one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------
$tx = $wb.Worksheets.Item("Transactions")

# The K/L columns only ever held empty, styled placeholder cells - drop them.
$tx.Columns("K:L").Delete()

# Row 2
$tx.Range("A2").Value = 6967
$tx.Range("B2").Value = "Head Office"
$tx.Range("C2").Value = 41713
$tx.Range("D2").Value = "Accrual"
$tx.Range("E2").NumberFormat = "General"
$tx.Range("E2").Value = 5.28
$tx.Range("F2").Value = 0
$tx.Range("G2").Value = 5.28
$tx.Range("H2").Value = 0
$tx.Range("I2").Value = 0
$tx.Range("J2").NumberFormat = "General"
$tx.Range("J2").Value = 0

# Row 3
$tx.Range("A3").Value = 6966
$tx.Range("B3").Value = "Head Office"
$tx.Range("C3").Value = 41708
$tx.Range("D3").Value = "Accrual"
$tx.Range("E3").NumberFormat = "General"
$tx.Range("E3").Value = 12.32
$tx.Range("F3").Value = 0
$tx.Range("G3").Value = 12.32
$tx.Range("H3").Value = 0
$tx.Range("I3").Value = 0
$tx.Range("J3").NumberFormat = "General"
$tx.Range("J3").Value = 0

# Row 4
$tx.Range("A4").Value = 6953
$tx.Range("B4").Value = "Head Office"
$tx.Range("C4").Value = 41708
$tx.Range("D4").Value = "Repayment"
$tx.Range("E4").NumberFormat = "#,##0"
$tx.Range("E4").Value = 1000
$tx.Range("F4").Value = 950.71
$tx.Range("G4").Value = 49.29
$tx.Range("H4").Value = 0
$tx.Range("I4").Value = 0
$tx.Range("J4").NumberFormat = "#,##0.00"
$tx.Range("J4").Value = 3214.17

# Row 5
$tx.Range("A5").Value = 6965
$tx.Range("B5").Value = "Head Office"
$tx.Range("C5").Value = 41699
$tx.Range("D5").Value = "Accrual"
$tx.Range("E5").NumberFormat = "General"
$tx.Range("E5").Value = 36.97
$tx.Range("F5").Value = 0
$tx.Range("G5").Value = 36.97
$tx.Range("H5").Value = 0
$tx.Range("I5").Value = 0
$tx.Range("J5").NumberFormat = "General"
$tx.Range("J5").Value = 0

# Row 6
$tx.Range("A6").Value = 6964
$tx.Range("B6").Value = "Head Office"
$tx.Range("C6").Value = 41672
$tx.Range("D6").Value = "Accrual"
$tx.Range("E6").NumberFormat = "General"
$tx.Range("E6").Value = 1.64
$tx.Range("F6").Value = 0
$tx.Range("G6").Value = 1.64
$tx.Range("H6").Value = 0
$tx.Range("I6").Value = 0
$tx.Range("J6").NumberFormat = "General"
$tx.Range("J6").Value = 0

# Row 7 (new)
$tx.Range("A7").Value = 6951
$tx.Range("B7").Value = "Head Office"
$tx.Range("C7").Value = 41672
$tx.Range("D7").Value = "Repayment"
$tx.Range("E7").NumberFormat = "General"
$tx.Range("E7").Value = 887.72
$tx.Range("F7").Value = 835.12
$tx.Range("G7").Value = 52.6
$tx.Range("H7").Value = 0
$tx.Range("I7").Value = 0
$tx.Range("J7").NumberFormat = "#,##0.00"
$tx.Range("J7").Value = 4164.88

# Row 8 (new)
$tx.Range("A8").Value = 6963
$tx.Range("B8").Value = "Head Office"
$tx.Range("C8").Value = 41671
$tx.Range("D8").Value = "Accrual"
$tx.Range("E8").NumberFormat = "General"
$tx.Range("E8").Value = 50.96
$tx.Range("F8").Value = 0
$tx.Range("G8").Value = 50.96
$tx.Range("H8").Value = 0
$tx.Range("I8").Value = 0
$tx.Range("J8").NumberFormat = "General"
$tx.Range("J8").Value = 0

# Row 9 (new)
$tx.Range("A9").Value = 6949
$tx.Range("B9").Value = "Head Office"
$tx.Range("C9").Value = 41640
$tx.Range("D9").Value = "Disbursement"
$tx.Range("E9").NumberFormat = "#,##0"
$tx.Range("E9").Value = 5000
$tx.Range("F9").Value = 0
$tx.Range("G9").Value = 0
$tx.Range("H9").Value = 0
$tx.Range("I9").Value = 0
$tx.Range("J9").NumberFormat = "#,##0"
$tx.Range("J9").Value = 5000

# Date columns keep their date format.
$tx.Range("C2:C9").NumberFormat = "d-mmm-yy"

# ---------------------------------------------------------------------
# DeleteAccountClosures sheet - drop the stray Navigate/LoanAccount row
# ---------------------------------------------------------------------
$dac = $wb.Worksheets.Item("DeleteAccountClosures")
$dac.Rows("4:4").Delete()

# ---------------------------------------------------------------------
# Repay1 sheet - selection cursor moved up one row
# ---------------------------------------------------------------------
$repay1 = $wb.Worksheets.Item("Repay1")
$repay1.Range("B2").Select()

# ---------------------------------------------------------------------
# Final selection / active sheet state
# ---------------------------------------------------------------------
$tx.Range("H9").Select()
$dac.Range("B2").Select()
